# Update generated ticket-demand figures for two events that appear on
# both the "展览" (exhibition) sheet and the "全部类型" (all types) sheet.
#   F2: 5331 -> 5359
#   F4: 920  -> 925

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5359
    $ws.Range("F4").Value = 925
}
